$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.763.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.286.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.10%  '
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.399'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.860.95'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('E13').Value = '  -3.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.839.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.312.08'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000163'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '435.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.437.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.511'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000113'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.193'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.29'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.11%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.64'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.39%  '
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.779.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.780'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0658'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '319.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.35%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0270'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.102'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
